# Apply the "fix: FHIR IG terminology and profile corrections" edit:
#  - Metadata!B7 (the "Experimental" row) gets the value "false" (was blank)
#  - Metadata!B8 (the "Date" row) is updated to the new generation timestamp
#
# Note: "false" must land in the sheet as literal TEXT (shared string),
# matching how the IG publisher tool authored it - NOT as an Excel boolean.
# Typing Value = "false" directly would auto-coerce to a Boolean (TRUE/FALSE)
# cell, so we enter it as text (leading apostrophe) and then copy the
# existing cell's format over it so the style index stays the same as the
# rest of the column (avoids leaving a stray "quote prefix" style behind).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Row 7: Experimental -> false ---
$expCell = $ws.Cells.Item(7, 2)
$expCell.Value = "'false"

# Re-apply the formatting of a neighboring plain-text cell (same original
# style) so the cell doesn't keep an Excel "stored as text" quote-prefix
# style that wasn't part of the original formatting.
$ws.Cells.Item(6, 2).Copy() | Out-Null
$expCell.PasteSpecial(-4122) | Out-Null

# --- Row 8: Date value update ---
$ws.Cells.Item(8, 2).Value = "2025-11-30T13:08:37+00:00"
